# Added Modem Configuration Page Test
#
# 1) Testdata sheet: add a new test-case row (row 9) for
#    "verifyUserAlreadyConnected", copying the row-8 formatting.
# 2) Locators sheet: add a new locator row (row 18) for the new
#    "ModemConfigurationPage" / "imAlreadyConnBtn" locator, copying
#    formatting from the appropriate existing cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Testdata sheet — new row 9
# ---------------------------------------------------------------
$testdata = $wb.Worksheets.Item("Testdata")

# Copy the formatting of the previous data row (row 8) down into row 9
$testdata.Range("A8:E8").Copy()
$testdata.Range("A9:E9").PasteSpecial(-4122)

$testdata.Range("A9").Value = "verifyUserAlreadyConnected"
$testdata.Range("B9").Value = "text1=Welcometo Globe myBusiness"
$testdata.Range("C9").Value = "mobileNumber=09271080510"
$testdata.Range("D9").Value = "pin=1111"
$testdata.Range("E9").Value = "nickname=Hendrix"

# ---------------------------------------------------------------
# Locators sheet — new row 18
# ---------------------------------------------------------------
$locators = $wb.Worksheets.Item("Locators")

# Column A (Page) — reuse the "text, vertical-bottom" box style already
# used for the IOS Locator Type column so the new row gets a closed box.
$locators.Range("E16").Copy()
$locators.Range("A18").PasteSpecial(-4122)
$locators.Range("A18").Value = "ModemConfigurationPage"

# Columns B:D (Locator Name / Android Locator Type / Android Locator) —
# reuse the existing "last row" box border, then force Text number format.
$locators.Range("B18").Copy()
$locators.Range("B18:D18").PasteSpecial(-4122)
$locators.Range("B18:D18").NumberFormat = "@"

$locators.Range("B18").Value = "imAlreadyConnBtn"
$locators.Range("C18").Value = "id"
$locators.Range("D18").Value = "ph.com.globe.mybusiness:id/btn_already_connected"
